# Fim da primeira semana
# Mark attendance "C" for columns K and L (29/04/2022 manhã/tarde)
# for every student row (3-49) that already has attendance marked in
# column J, mirroring the existing pattern used for columns C-J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 49; $row++) {
    $jValue = $ws.Cells.Item($row, 10).Value2
    if ($jValue -eq "C") {
        $ws.Cells.Item($row, 11).Value = "C"
        $ws.Cells.Item($row, 12).Value = "C"
    }
}

$ws.Range("L3").Select()
